# Append new scraped rows (2025-12-06 12:34:19 JST run) to the "ランサーズ" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- column widths (B: 35 -> 51, H: 13 -> 16) -------------------------------
# ColumnWidth setter adds Excel's default padding (~0.8333 chars) before it is
# stored as the <col width="..."> attribute, so back that padding out here in
# order to land on the exact target widths.
$pad = 5 / 6
$ws.Columns.Item(2).ColumnWidth = 51 - $pad
$ws.Columns.Item(8).ColumnWidth = 16 - $pad

# --- row data ----------------------------------------------------------------
$timestamp = "2025-12-06 12:34:19"

$rows = @(
    @($timestamp, "CapcutAPIを用いた動画の自動制作ツールの作成", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448662", 248, "🔥API ◆ツール"),
    @($timestamp, "CapcutAPIを用いた動画の自動制作ツールの作成", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448659", 248, "🔥API ◆ツール"),
    @($timestamp, "【自動化】Webサービス更新ツール開発(200アカウント管理)", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448409", 230, "◆ツール,開発 ◇管理"),
    @($timestamp, "【品質重視】出張買取サービス向け予約管理システム開発(UI/要件定義済/Cursor実装途中あり)", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448677", 153, "◆開発,システム開発 ◇管理"),
    @($timestamp, "【急募】新規システム開発に伴う要件定義依頼", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448563", 110, "◆開発,システム開発"),
    @($timestamp, "【受注メールを元にECサイト自動仕入ツール】", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448428", 98, "◆ツール ◇サイト"),
    @($timestamp, "【緊急】既存コードの構造解析ができるPHPエンジニアを探しています", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448440", 33, "○PHP")
)

# drop the old hyperlink metadata; it will be rebuilt below with the refreshed
# row positions (values/formatting for the cells themselves stay intact)
$ws.Hyperlinks.Delete()

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]

    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row[5])
    $ws.Cells.Item($r, 6).Style = "Hyperlink"

    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
